$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 201
$ws.Range("I2").Value = 521
$ws.Range("J2").Value = 2219
$ws.Range("K2").Value = 9
$ws.Range("L2").Value = 561
$ws.Range("M2").Value = 36
$ws.Range("N2").Value = 391
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 11
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 32
$ws.Range("S2").Value = 244
$ws.Range("T2").Value = 403
$ws.Range("U2").Value = 28
$ws.Range("V2").Value = 3341
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 3401
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 51
$ws.Range("AA2").Value = 19
